$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (nRF52832): Flash (kB) = "<=512"
$ws.Range("K5").Value = "<=512"

# Row 6 (nRF52833): Flash (kB) = 512
$ws.Range("K6").Value = 512

# Row 7 (nRF52840): Flash (kB) = 1024, Zusatz = module list, row grows to 3 lines
$ws.Range("K7").Value = 1024
$ws.Range("M7").Value = "ANT: NINA-B302-00B (u-blox)`n          BMD-340-A-R (u-blox)`n          BL654PA (Lairdconnect)"
$ws.Range("M7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 43.2

# Row 8 (nRF5340): Flash (kB) = 1024, Zusatz = module list
$ws.Range("K8").Value = 1024
$ws.Range("M8").Value = "ANT: NORA-B106-00B (u-blox)`n          BL5340 (Lairdconnect)"
$ws.Range("M8").WrapText = $true

# Row 10 (EFR32BG27): Zusatz = "-"
$ws.Range("M10").Value = "-"

# Row 12 (BlueNRG-LP): Zusatz = "-"
$ws.Range("M12").Value = "-"

# Column M width adjustment (bestFit column got a bit wider for the new text)
$ws.Columns.Item(13).ColumnWidth = 26.75

# Final selection as left by the author
$ws.Range("M11").Select()
